$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The sheet grows from 15 to 20 rows: a new "Business Contact" block
# (4 rows) is inserted, plus one extra "Department / Agency" detail
# row is added to both the "Business Sponsor" and "Technical Contact"
# blocks. We rebuild rows 6-20 explicitly (bottom-up for the copies,
# so a source row is always read before it is overwritten), then fix
# up every cell's text/blankness explicitly afterwards (Copy() here
# only reliably transfers *formatting*, not values).
# ------------------------------------------------------------------

# --- Shift the "Technical Contact" block (old rows 9-15) down to new rows 14-20 ---
$ws.Range("A15:B15").Copy($ws.Range("A20:B20"))
$ws.Range("A14:B14").Copy($ws.Range("A19:B19"))
$ws.Range("A13:B13").Copy($ws.Range("A18:B18"))
$ws.Range("A12:B12").Copy($ws.Range("A17:B17"))
$ws.Range("A11:B11").Copy($ws.Range("A16:B16"))
$ws.Range("A10:B10").Copy($ws.Range("A15:B15"))
$ws.Range("A9:B9").Copy($ws.Range("A14:B14"))

# --- Build the brand-new "Business Contact" block (new rows 10-13) from the still-pristine "Business Sponsor" rows 6-8 ---
$ws.Range("A6:B6").Copy($ws.Range("A10:B10"))
$ws.Range("A7:B7").Copy($ws.Range("A11:B11"))
$ws.Range("A8:B8").Copy($ws.Range("A12:B12"))
$ws.Range("A8:B8").Copy($ws.Range("A13:B13"))

# --- Add the new "Department / Agency" row (new row 9) to the Business Sponsor block ---
$ws.Range("A8:B8").Copy($ws.Range("A9:B9"))

# ------------------------------------------------------------------
# Now make every cell's value exactly right (Copy only reliably
# carries over formatting - not values - in this environment, so
# every destination cell above is explicitly corrected below).
# ------------------------------------------------------------------

$ws.Range("A1").Value = "`nNSLS Certificate Signing Request Form"
$ws.Range("B1").ClearContents()

$ws.Range("A2").ClearContents()
$ws.Range("B2").Value = "Application Information"

$ws.Range("A3").Value = "Dept. / Agency Name  "
$ws.Range("B3").ClearContents()

$ws.Range("A4").Value = "Application Name  "
$ws.Range("B4").ClearContents()

$ws.Range("A5").Value = "Environment [-select one-]  "
$ws.Range("B5").ClearContents()

$ws.Range("A6").ClearContents()
$ws.Range("B6").Value = "Business Sponsor (GNS/NSHA)"

$ws.Range("A7").Value = "Name  "
$ws.Range("B7").ClearContents()

$ws.Range("A8").Value = "Email  "
$ws.Range("B8").ClearContents()

$ws.Range("A9").Value = "Department  / Agency"
$ws.Range("B9").ClearContents()

$ws.Range("A10").ClearContents()
$ws.Range("B10").Value = "Business Contact"

$ws.Range("A11").Value = "Name  "
$ws.Range("B11").ClearContents()

$ws.Range("A12").Value = "Email  "
$ws.Range("B12").ClearContents()

$ws.Range("A13").Value = "Department / Agency /Vendor Name  "
$ws.Range("B13").ClearContents()

$ws.Range("A14").ClearContents()
$ws.Range("B14").Value = "Technical Contact"

$ws.Range("A15").Value = "Name  "
$ws.Range("B15").ClearContents()

$ws.Range("A16").Value = "Email  "
$ws.Range("B16").ClearContents()

$ws.Range("A17").Value = "Phone  "
$ws.Range("B17").ClearContents()

$ws.Range("A18").Value = "Department / Agency /Vendor Name  "
$ws.Range("B18").ClearContents()

# Row 19 is a bare blank spacer row - never touch A19/B19 cell content, just its height (set below).

# Row 20 keeps the rich-text "Note" - untouched text, but make sure no stray content leaks in.
$ws.Range("B20").ClearContents()

# ------------------------------------------------------------------
# Row heights for the rows that are new past the original dimension.
# ------------------------------------------------------------------
$ws.Range("A9").RowHeight = 21
$ws.Range("A10").RowHeight = 21
$ws.Range("A11").RowHeight = 21
$ws.Range("A12").RowHeight = 21
$ws.Range("A13").RowHeight = 21
$ws.Range("A14").RowHeight = 21
$ws.Range("A15").RowHeight = 21
$ws.Range("A16").RowHeight = 21
$ws.Range("A17").RowHeight = 21
$ws.Range("A18").RowHeight = 21
$ws.Range("A19").RowHeight = 21.75
$ws.Range("A20").RowHeight = 21

# ------------------------------------------------------------------
# Merged cells: the title merge (A1:B1) is untouched; the "Note" merge
# moves from the old A15:B15 down to the new A20:B20.
# ------------------------------------------------------------------
if ($ws.Range("A15:B15").MergeCells) { $ws.Range("A15:B15").UnMerge() }
$ws.Range("A20:B20").Merge() | Out-Null

# ------------------------------------------------------------------
# Column A: width 40 (stored width "40" once Excel's 0.8333 padding
# is accounted for), keep custom width.
# ------------------------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 39.1666666666667

# ------------------------------------------------------------------
# Sheet view selection: now a single cell, B3.
# ------------------------------------------------------------------
$ws.Range("B3").Select()
